# Apply weekly update: a new data row is inserted before current row 126
# (shifting all following rows down by one) and populated with the
# latest weekly price record for Jengibre / Primera.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 126, pushing existing rows 126-131 down to 127-132.
$ws.Rows.Item(126).Insert()

# Populate the newly inserted row 126 with the new weekly record.
$ws.Cells.Item(126, 1).Value = 9
$ws.Cells.Item(126, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(126, 3).Value = "Metropolitana"
$ws.Cells.Item(126, 4).Value = 45075
$ws.Cells.Item(126, 5).Value = 13
$ws.Cells.Item(126, 6).Value = 100114007
$ws.Cells.Item(126, 7).Value = "Jengibre"
$ws.Cells.Item(126, 8).Value = "Sin especificar"
$ws.Cells.Item(126, 9).Value = "Primera"
$ws.Cells.Item(126, 10).Value = 520
$ws.Cells.Item(126, 11).Value = 16000
$ws.Cells.Item(126, 12).Value = 17000
$ws.Cells.Item(126, 13).Value = 16500
$ws.Cells.Item(126, 14).Value = "$/caja 13 kilos"
$ws.Cells.Item(126, 15).Value = "Perú"
$ws.Cells.Item(126, 16).Value = 1269
$ws.Cells.Item(126, 17).Value = 13
$ws.Cells.Item(126, 18).Value = "Hortaliza"
